$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the surviving action-card rows with the re-imported category names
$ws.Range("B41").Value = "action card 5"
$ws.Range("B42").Value = "Emotinal health"
$ws.Range("B43").Value = "Happyness"
$ws.Range("B44").Value = "safe-delivery"

# Remove the now-obsolete trailing rows (old rows 45-49)
$ws.Range("A45:B49").EntireRow.Delete()
